# ==========================================================================
# Edit script for LOBSTAHS_adduct_ion_hierarchies.xlsx
# Implements:
#  1) Sheet "Adduct ion hierarchies": extend table with 15 new lipid-class
#     columns (Q:AE), extend/re-merge the title row, fill in new header and
#     data cells, select C1:AE1 and make this sheet the active tab.
#  2) Sheet "Notes": rewrite the notes/history content, add 4 new file
#     history rows, resize columns A and B, select A22, and make it the
#     non-active tab.
# ==========================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Adduct ion hierarchies")
$ws2 = $wb.Worksheets.Item("Notes")

# --------------------------------------------------------------------------
# 1) "Adduct ion hierarchies" sheet
# --------------------------------------------------------------------------

# Extend the merged/styled title band from C1:P1 to C1:AE1
$ws1.Range("C1:P1").UnMerge() | Out-Null
$ws1.Range("C1").Copy() | Out-Null
$ws1.Range("Q1:AE1").PasteSpecial(-4122) | Out-Null
$ws1.Range("C1:AE1").Merge() | Out-Null

# New header labels for the 15 added lipid classes
$ws1.Range("Q2").Value  = "vGSL"
$ws1.Range("R2").Value  = "sGSL"
$ws1.Range("S2").Value  = "hGSL"
$ws1.Range("T2").Value  = "hapGSL"
$ws1.Range("U2").Value  = "PDPT"
$ws1.Range("V2").Value  = "BLL"
$ws1.Range("W2").Value  = "hapCER"
$ws1.Range("X2").Value  = "S_DGCC"
$ws1.Range("Y2").Value  = "fungalGSL"
$ws1.Range("Z2").Value  = "plastoquinone_9"
$ws1.Range("AA2").Value = "plastoquinone_9OH"
$ws1.Range("AB2").Value = "plastoquinone_9OH2"
$ws1.Range("AC2").Value = "ubiquinone"
$ws1.Range("AD2").Value = "scytonemin"
$ws1.Range("AE2").Value = "reduced_scytonemin"

# Abundance-rank data for the new columns (rows 3-22)
$ws1.Range("Q3").Value = 5
$ws1.Range("R3").Value = 1
$ws1.Range("S3").Value = 2
$ws1.Range("U3").Value = 4
$ws1.Range("V3").Value = 1
$ws1.Range("W3").Value = 3
$ws1.Range("Y3").Value = 2
$ws1.Range("Z3").Value = 1
$ws1.Range("AC3").Value = 1
$ws1.Range("AD3").Value = 1
$ws1.Range("V4").Value = 2
$ws1.Range("Q5").Value = 3
$ws1.Range("S5").Value = 3
$ws1.Range("T5").Value = 3
$ws1.Range("W5").Value = 2
$ws1.Range("X5").Value = 3
$ws1.Range("Q7").Value = 1
$ws1.Range("S7").Value = 1
$ws1.Range("T7").Value = 1
$ws1.Range("U7").Value = 1
$ws1.Range("W7").Value = 1
$ws1.Range("X7").Value = 1
$ws1.Range("Y7").Value = 1
$ws1.Range("Q8").Value = 8
$ws1.Range("R8").Value = 2
$ws1.Range("S8").Value = 5
$ws1.Range("Q9").Value = 9
$ws1.Range("S9").Value = 7
$ws1.Range("Q10").Value = 6
$ws1.Range("T10").Value = 5
$ws1.Range("U10").Value = 5
$ws1.Range("X10").Value = 5
$ws1.Range("Q11").Value = 2
$ws1.Range("T11").Value = 2
$ws1.Range("U11").Value = 2
$ws1.Range("W11").Value = 4
$ws1.Range("X11").Value = 2
$ws1.Range("Y11").Value = 3
$ws1.Range("Q12").Value = 7
$ws1.Range("S12").Value = 6
$ws1.Range("T12").Value = 6
$ws1.Range("U12").Value = 6
$ws1.Range("X12").Value = 6
$ws1.Range("Q13").Value = 4
$ws1.Range("S13").Value = 4
$ws1.Range("T13").Value = 4
$ws1.Range("U13").Value = 3
$ws1.Range("W13").Value = 5
$ws1.Range("X13").Value = 4
$ws1.Range("Y13").Value = 4
$ws1.Range("Q14").Value = 1
$ws1.Range("R14").Value = 1
$ws1.Range("S14").Value = 1
$ws1.Range("T14").Value = 3
$ws1.Range("U14").Value = 1
$ws1.Range("V14").Value = 1
$ws1.Range("W14").Value = 1
$ws1.Range("X14").Value = 1
$ws1.Range("Y14").Value = 1
$ws1.Range("Z14").Value = 1
$ws1.Range("AA14").Value = 3
$ws1.Range("AB14").Value = 6
$ws1.Range("AC14").Value = 4
$ws1.Range("AD14").Value = 1
$ws1.Range("AE14").Value = 1
$ws1.Range("Q15").Value = 7
$ws1.Range("R15").Value = 7
$ws1.Range("S15").Value = 5
$ws1.Range("T15").Value = 6
$ws1.Range("U15").Value = 5
$ws1.Range("X15").Value = 4
$ws1.Range("Y15").Value = 4
$ws1.Range("Z15").Value = 5
$ws1.Range("AA15").Value = 6
$ws1.Range("AB15").Value = 4
$ws1.Range("AC15").Value = 6
$ws1.Range("Q16").Value = 3
$ws1.Range("R16").Value = 4
$ws1.Range("S16").Value = 3
$ws1.Range("T16").Value = 4
$ws1.Range("Z16").Value = 2
$ws1.Range("AA16").Value = 1
$ws1.Range("AB16").Value = 1
$ws1.Range("AC16").Value = 1
$ws1.Range("Q17").Value = 2
$ws1.Range("R17").Value = 2
$ws1.Range("S17").Value = 2
$ws1.Range("T17").Value = 1
$ws1.Range("U17").Value = 3
$ws1.Range("W17").Value = 2
$ws1.Range("X17").Value = 2
$ws1.Range("Y17").Value = 2
$ws1.Range("Z17").Value = 4
$ws1.Range("AA17").Value = 2
$ws1.Range("AB17").Value = 2
$ws1.Range("AC17").Value = 3
$ws1.Range("AD17").Value = 2
$ws1.Range("Q18").Value = 8
$ws1.Range("R18").Value = 3
$ws1.Range("T18").Value = 7
$ws1.Range("U18").Value = 6
$ws1.Range("AB18").Value = 9
$ws1.Range("Q19").Value = 5
$ws1.Range("R19").Value = 5
$ws1.Range("S19").Value = 4
$ws1.Range("T19").Value = 2
$ws1.Range("U19").Value = 2
$ws1.Range("V19").Value = 2
$ws1.Range("W19").Value = 4
$ws1.Range("Y19").Value = 5
$ws1.Range("Z19").Value = 3
$ws1.Range("AA19").Value = 4
$ws1.Range("AB19").Value = 3
$ws1.Range("AC19").Value = 2
$ws1.Range("Q20").Value = 6
$ws1.Range("R20").Value = 6
$ws1.Range("S20").Value = 6
$ws1.Range("T20").Value = 5
$ws1.Range("U20").Value = 7
$ws1.Range("W20").Value = 5
$ws1.Range("X20").Value = 5
$ws1.Range("AA20").Value = 5
$ws1.Range("AB20").Value = 5
$ws1.Range("AC20").Value = 7
$ws1.Range("AB21").Value = 8
$ws1.Range("Q22").Value = 4
$ws1.Range("S22").Value = 7
$ws1.Range("U22").Value = 4
$ws1.Range("W22").Value = 3
$ws1.Range("X22").Value = 3
$ws1.Range("Y22").Value = 3
$ws1.Range("AB22").Value = 7
$ws1.Range("AC22").Value = 5

# --------------------------------------------------------------------------
# 2) "Notes" sheet
# --------------------------------------------------------------------------

# Insert one row above the old row 2, shifting everything else down by one
$ws2.Rows.Item(2).Insert() | Out-Null

# New row 2 / row 3: source-data citation (split across two rows/cells)
$ws2.Range("A3").ClearContents() | Out-Null
$ws2.Range("A2").Value = "Source data:"
$ws2.Range("B2").Value = "Table 2 in Collins, J. R., B. R. Edwards, H. F. Fredricks, and B. A. S. Van Mooy. 2016. LOBSTAHS: An adduct-based lipidomics strategy for discovery and identification of oxidative stress biomarkers. Analytical Chemistry 88: 7154-7162; doi:10.1021/acs.analchem.6b01260"
$ws2.Range("B3").Value = "For BLL, PDPT, vGSL, sGSL, hGSL, hapGSL, and hapCER: Hunter J. E., M. J. Frada, H. F. Fredricks, A. Vardi, and B. A. S. Van Mooy. 2015. Targeted and untargeted lipidomics of Emiliania huxleyi viral infection and life cycle phases highlights molecular biomarkers of infection, susceptibility, and ploidy. Front. Mar. Sci. 2: 81; doi: 10.3389/fmars.2015.00081"

# Row 4 (retains the black-font style that was on the old row 3) gets new wording
$ws2.Range("A4").Value = 'The first tab of this workbook can be used to generate the file "LOBSTAHS_adduct_ion_hierarchies.csv," required for lipid-oxlipid-oxyipin database generation in the LOBSTAHS lipidomics screening pipeline'

# Rows 5, 6, 8, 10, 12, 13 keep their old (shifted) text already in place -
# no edits required there.

# Append four new file-history rows (18-21), reusing the date style (s=2)
# from the existing history rows by copy/paste, then filling in the values
$ws2.Range("A17").Copy() | Out-Null
$ws2.Range("A18:A21").PasteSpecial(-4122) | Out-Null

$ws2.Range("A18").Value = 42624
$ws2.Range("B18").Value = "Added BLL, PDPT, vGSL, sGSL, hGSL, hapGSL, and hapCER (from Hunter et al., 2015; citation above) "
$ws2.Range("C18").Value = "JEH"

$ws2.Range("A19").Value = 42709
$ws2.Range("B19").Value = "Added S_DGCC, fungalGSLs, ubiquinones and plastoquinones"
$ws2.Range("C19").Value = "JEH"

$ws2.Range("A20").Value = 42719
$ws2.Range("B20").Value = "Added scytonemin and reduced scytonemin"
$ws2.Range("C20").Value = "JEH"

$ws2.Range("A21").Value = 42758
$ws2.Range("B21").Value = "Modifications as necessary for R function upgrades"
$ws2.Range("C21").Value = "JRC"

# Resize columns A and B (target raw widths 12.5 / 74.5 characters); Excel's
# ColumnWidth setter adds ~0.8333 of internal padding, so compensate for it
$ws2.Columns.Item(1).ColumnWidth = 12.5 - 0.8333333333333333
$ws2.Columns.Item(2).ColumnWidth = 74.5 - 0.8333333333333333

# --------------------------------------------------------------------------
# 3) Selections / active tab
# --------------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("A22").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("C1:AE1").Select() | Out-Null
